{"js": "// Update the date heading and every two-digit x two-digit multiplication\n// answer in the practice-sheet table. Each \"before\" string is unique in\n// the document, so a targeted search + replace is safe and preserves the\n// existing run formatting (font/size) on every cell.\nconst pairs = [\n  ['2025-03-18 Tuesday', '2025-03-19 Wednesday'],\n  ['12\u00d774=888', '12\u00d782=984'],\n  ['49\u00d747=2303', '56\u00d757=3192'],\n  ['88\u00d794=8272', '46\u00d784=3864'],\n  ['88\u00d723=2024', '56\u00d750=2800'],\n  ['51\u00d735=1785', '47\u00d763=2961'],\n  ['38\u00d752=1976', '75\u00d786=6450'],\n  ['75\u00d792=6900', '89\u00d731=2759'],\n  ['49\u00d757=2793', '51\u00d742=2142'],\n  ['79\u00d740=3160', '63\u00d773=4599'],\n  ['60\u00d796=5760', '94\u00d759=5546'],\n  ['58\u00d744=2552', '39\u00d718=702'],\n  ['45\u00d784=3780', '83\u00d790=7470'],\n  ['60\u00d754=3240', '92\u00d796=8832'],\n  ['85\u00d774=6290', '34\u00d763=2142'],\n  ['33\u00d780=2640', '22\u00d740=880'],\n  ['19\u00d733=627', '72\u00d718=1296'],\n  ['16\u00d763=1008', '50\u00d725=1250'],\n  ['70\u00d756=3920', '80\u00d776=6080'],\n  ['89\u00d739=3471', '16\u00d738=608'],\n  ['64\u00d792=5888', '49\u00d723=1127'],\n  ['98\u00d718=1764', '61\u00d723=1403'],\n  ['20\u00d716=320', '67\u00d740=2680'],\n  ['19\u00d762=1178', '34\u00d764=2176'],\n  ['91\u00d775=6825', '63\u00d787=5481'],\n  ['84\u00d777=6468', '32\u00d738=1216'],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n\n  // Replace in-place so the run's existing formatting (rFonts/sz) sticks.\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Update the date heading and every two-digit x two-digit multiplication\n# answer in the practice-sheet table. Each \"before\" string is unique in\n# the document, so Find/Replace on the whole-document range is safe and\n# leaves the surrounding run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-03-18 Tuesday\", \"2025-03-19 Wednesday\"),\n    @(\"12\u00d774=888\", \"12\u00d782=984\"),\n    @(\"49\u00d747=2303\", \"56\u00d757=3192\"),\n    @(\"88\u00d794=8272\", \"46\u00d784=3864\"),\n    @(\"88\u00d723=2024\", \"56\u00d750=2800\"),\n    @(\"51\u00d735=1785\", \"47\u00d763=2961\"),\n    @(\"38\u00d752=1976\", \"75\u00d786=6450\"),\n    @(\"75\u00d792=6900\", \"89\u00d731=2759\"),\n    @(\"49\u00d757=2793\", \"51\u00d742=2142\"),\n    @(\"79\u00d740=3160\", \"63\u00d773=4599\"),\n    @(\"60\u00d796=5760\", \"94\u00d759=5546\"),\n    @(\"58\u00d744=2552\", \"39\u00d718=702\"),\n    @(\"45\u00d784=3780\", \"83\u00d790=7470\"),\n    @(\"60\u00d754=3240\", \"92\u00d796=8832\"),\n    @(\"85\u00d774=6290\", \"34\u00d763=2142\"),\n    @(\"33\u00d780=2640\", \"22\u00d740=880\"),\n    @(\"19\u00d733=627\", \"72\u00d718=1296\"),\n    @(\"16\u00d763=1008\", \"50\u00d725=1250\"),\n    @(\"70\u00d756=3920\", \"80\u00d776=6080\"),\n    @(\"89\u00d739=3471\", \"16\u00d738=608\"),\n    @(\"64\u00d792=5888\", \"49\u00d723=1127\"),\n    @(\"98\u00d718=1764\", \"61\u00d723=1403\"),\n    @(\"20\u00d716=320\", \"67\u00d740=2680\"),\n    @(\"19\u00d762=1178\", \"34\u00d764=2176\"),\n    @(\"91\u00d775=6825\", \"63\u00d787=5481\"),\n    @(\"84\u00d777=6468\", \"32\u00d738=1216\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    # wdFindContinue=1, wdReplaceOne=1 -> replace the single (unique) match.\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
